# Applies the "added harvard case classification" update:
# recomputes the "_old" evaluation columns (Ada_old, Avey_old, WebMD_old,
# doctor_NJ_old, doctor_TH_old) for the precision/recall/f1/f2 rows, and
# the corresponding "length (x of gs)" row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (recall)
$ws.Range("C3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("Q3").Value = 1
$ws.Range("U3").Value = 1
$ws.Range("W3").Value = 1

# Row 4 (f1-score)
$ws.Range("C4").Value = 0.3333333333333334
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("Q4").Value = 0.4
$ws.Range("U4").Value = 0.6666666666666666
$ws.Range("W4").Value = 0.5

# Row 5 (f2-score)
$ws.Range("C5").Value = 0.5555555555555556
$ws.Range("F5").Value = 0.8333333333333334
$ws.Range("Q5").Value = 0.625
$ws.Range("U5").Value = 0.8333333333333334
$ws.Range("W5").Value = 0.7142857142857143

# Row 11 (length (x of gs))
$ws.Range("C11").Value = 5
$ws.Range("F11").Value = 2
$ws.Range("I11").Value = 2
$ws.Range("M11").Value = 2
$ws.Range("Q11").Value = 4
$ws.Range("S11").Value = 1
$ws.Range("U11").Value = 2
$ws.Range("W11").Value = 3
